$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 308
$wsExhibit.Range("F4").Value = 214
$wsExhibit.Range("F5").Value = 2630
$wsExhibit.Range("F6").Value = 1879
$wsExhibit.Range("F7").Value = 360
$wsExhibit.Range("F9").Value = 931

# Sheet "演出" (performances) - update column F
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 26

# Sheet "全部类型" (all types) - update column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 308
$wsAll.Range("F4").Value = 214
$wsAll.Range("F5").Value = 2630
$wsAll.Range("F6").Value = 1879
$wsAll.Range("F7").Value = 360
$wsAll.Range("F8").Value = 26
$wsAll.Range("F10").Value = 931
